# Append one new data row to each of the four worksheets, mirroring the
# structure of the existing rows (time / hex-string fields / decimal fields).

$wb = $excel.ActiveWorkbook

$dateFormat = "YYYY-MM-DD HH:MM:SS"
$bigNumber = 568631262647113970876416.0

function Add-DataRow {
    param(
        $ws,
        $row,
        $timeValue,
        $colB,
        $colC,
        $colD,
        $colE,
        $colF,
        $colG,
        $colH,
        $colI,
        $colGIsText
    )

    $ws.Cells.Item($row, 1).Value = $timeValue
    $ws.Cells.Item($row, 1).NumberFormat = $dateFormat

    $ws.Cells.Item($row, 2).Value = $colB
    $ws.Cells.Item($row, 3).Value = $colC
    $ws.Cells.Item($row, 4).Value = $colD
    $ws.Cells.Item($row, 5).Value = $colE

    $ws.Cells.Item($row, 6).Value = $colF

    if ($colGIsText) {
        # This particular value is too large to round-trip exactly through a
        # double, so the source data stores it as literal text instead of a
        # number. Force text storage, then strip the formatting footprint so
        # the cell ends up unstyled, matching the rest of the row.
        $ws.Cells.Item($row, 7).NumberFormat = "@"
        $ws.Cells.Item($row, 7).Value = $colG
        $ws.Cells.Item($row, 7).ClearFormats()
    } else {
        $ws.Cells.Item($row, 7).Value = $colG
    }

    $ws.Cells.Item($row, 8).Value = $colH
    $ws.Cells.Item($row, 9).Value = $colI
}

# ROW50-FE-LIFTER (sheet 1): new row 58
$ws1 = $wb.Worksheets.Item(1)
Add-DataRow $ws1 58 45753.21278416667 "0x01,0x90" "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x15,0x41,0x0c," "0x01,0x5e" "0xe" 400 $bigNumber 350 14 $false

# ROW50-MID-LIFTER (sheet 2): new row 60
$ws2 = $wb.Worksheets.Item(2)
Add-DataRow $ws2 60 45753.18106481482 "0x01,0x90 " "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x20,0x41,0x0c," "0x01,0x62" "0x19" 400 "568631262647113771663628" 354 25 $true

# ROW11-FE-LIFTER (sheet 3): new row 58
$ws3 = $wb.Worksheets.Item(3)
Add-DataRow $ws3 58 45753.24469516204 "0x01,0x90" "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x1a,0x41,0x0c," "0x01,0x5e" "0x14" 400 $bigNumber 350 20 $false

# ROW11-MID-LIFTER (sheet 4): new row 58
$ws4 = $wb.Worksheets.Item(4)
Add-DataRow $ws4 58 45753.37799475694 "0x01,0x90" "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x11,0x40,0x0c," "0x01,0x62" "0x19" 400 $bigNumber 354 25 $false
